$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer present in the final table:
#  row 22 -> code "EF" ("?")
#  row 21 -> code "EA" (blank name/values)
#  row 17 -> code "AC" ("PIP Video Input")
#  row 13 -> code "72" ("Gamma")
# Deleted from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(13).Delete()

# Expand the "Standard Color" possible values (row 4, code "14") with the full list of calibration modes.
$standardColor = "1 - sRGB`n2 - `n4 - `n5 - `n6 - `n8 - `n11 - `n14 - Adobe`n15 - EBU`n16 - SMPTE-C`n17 - REC709`n18 - DICOM SIM`n19 - DCI-P3`n21 - CAL 1`n22 - CAL 2`n23 - CAL 3`n24 - iPhone (DCI-P3)`n255 - Custom"
$ws.Cells.Item(4, 4).Value = $standardColor

# Expand the "Display Application" possible values (now row 16, code "DC") with the full list of applications.
$displayApplication = "0 - Default`n3 - Movie`n4 - ?`n48 - Game FPS1`n49 - Game FPS2`n50 - Game RTS`n51 - Game MODA`n52 - Web`n53 - Text`n54 - MAC`n55 - Designer CAD/CAM`n56 - Designer Animation`n57 - Designer VideoEdit`n58 - Photographer Retro`n59 - Photographer Photo`n60 - Photographer Landscape`n61 - Photographer Portrait`n62 - Photographer Monochrome"
$ws.Cells.Item(16, 4).Value = $displayApplication

# Restore the row heights for the two cells that now hold much longer text.
$ws.Rows.Item(4).RowHeight = 261
$ws.Rows.Item(16).RowHeight = 261

# Update the view: select D4 (matching the saved view state).
$ws.Range("D4").Select()
